$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 84.72727
$ws.Range("I5").Value = 84.72727
$ws.Range("K5").Value = 84.72727
$ws.Range("M5").Value = 30.27273

$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = ""

$ws.Range("H132").Value = 5291.3335
$ws.Range("I132").Value = 5291.3335
$ws.Range("K132").Value = 15874.0005
$ws.Range("M132").Value = -13344.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1864.6
$ws.Range("J2").Value = 913
$ws.Range("L2").Value = 913
$ws.Range("N2").Value = -1139

$ws.Range("H45").Value = 15998.75
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 15998.75
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15998.75
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -16752.75

$ws.Range("H61").Value = 9877.272000000001
$ws.Range("I61").Value = 7956.25
$ws.Range("K61").Value = 7956.25
$ws.Range("M61").Value = -7744.25

$ws.Range("H74").Value = 7476.875
$ws.Range("I74").Value = 3302.5
$ws.Range("K74").Value = 3302.5
$ws.Range("M74").Value = -2428.5

$ws.Range("H77").Value = 7476.875
$ws.Range("I77").Value = 3302.5
$ws.Range("K77").Value = 16512.5
$ws.Range("M77").Value = -12144.5

$ws.Range("H116").Value = 1864.6
$ws.Range("J116").Value = 913
$ws.Range("L116").Value = 913
$ws.Range("N116").Value = -5501

$ws.Range("H132").Value = 12082.833
$ws.Range("I132").Value = 8166
$ws.Range("K132").Value = 24498
$ws.Range("M132").Value = -21968

$ws.Range("H136").Value = 9877.272000000001
$ws.Range("I136").Value = 7956.25
$ws.Range("K136").Value = 23868.75
$ws.Range("M136").Value = -21318.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1864.6
$ws.Range("J3").Value = 913
$ws.Range("L3").Value = 913
$ws.Range("N3").Value = -1141

$ws.Range("H35").Value = 38699.168
$ws.Range("J35").Value = 38699.168
$ws.Range("L35").Value = 38699.168
$ws.Range("N35").Value = -39319.168

$ws.Range("H107").Value = 979.4
$ws.Range("I107").Value = 966.3333
$ws.Range("K107").Value = 966.3333
$ws.Range("M107").Value = 953.6667

$ws.Range("H126").Value = 20000
$ws.Range("J126").Value = 20000
$ws.Range("L126").Value = 20000
$ws.Range("N126").Value = -29880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 238.66667
$ws.Range("I15").Value = 238.66667
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 238.66667
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -68.66667000000001
$ws.Range("N15").Value = ""

$ws.Range("H19").Value = 233.33333
$ws.Range("I19").Value = 160
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 160
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = 10
$ws.Range("N19").Value = -940

$ws.Range("H24").Value = 233.33333
$ws.Range("I24").Value = 160
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 160
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 10
$ws.Range("N24").Value = -940

$ws.Range("H82").Value = 64999
$ws.Range("J82").Value = 64999
$ws.Range("L82").Value = 64999
$ws.Range("N82").Value = -65721

$ws.Range("H85").Value = 64999
$ws.Range("J85").Value = 64999
$ws.Range("L85").Value = 64999
$ws.Range("N85").Value = -67495

$ws.Range("H132").Value = 2738.7727
$ws.Range("I132").Value = 1515.0588
$ws.Range("K132").Value = 4545.1764
$ws.Range("M132").Value = -2015.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 9.199999999999999
$ws.Range("J7").Value = 10
$ws.Range("L7").Value = 30
$ws.Range("N7").Value = -254

$ws.Range("H13").Value = 1322.2
$ws.Range("I13").Value = 1320.3334
$ws.Range("K13").Value = 3961.0002
$ws.Range("M13").Value = -3793.0002

$ws.Range("H116").Value = 1062.75
$ws.Range("I116").Value = 1062.75
$ws.Range("K116").Value = 3188.25
$ws.Range("M116").Value = 253.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 16.705883
$ws.Range("I2").Value = 16.555555
$ws.Range("J2").Value = 16.875
$ws.Range("K2").Value = 16.555555
$ws.Range("L2").Value = 16.875
$ws.Range("M2").Value = 96.444445
$ws.Range("N2").Value = -242.875

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H118").Value = 23435.666
$ws.Range("J118").Value = 23435.666
$ws.Range("L118").Value = 23435.666
$ws.Range("N118").Value = -26749.666

$ws.Range("H122").Value = 1099.5
$ws.Range("I122").Value = 1099.5
$ws.Range("K122").Value = 3298.5
$ws.Range("M122").Value = -848.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1125
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = ""

$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622

$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4833.1665
$ws.Range("J14").Value = 4833.1665
$ws.Range("L14").Value = 4833.1665
$ws.Range("N14").Value = -5169.1665

$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812

$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808

$ws.Range("H132").Value = 10922.846
$ws.Range("I132").Value = 8249.5
$ws.Range("J132").Value = 12111
$ws.Range("K132").Value = 24748.5
$ws.Range("L132").Value = 36333
$ws.Range("M132").Value = -22218.5

Write-Output "applied all changes"
